$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'net cycling shorts',
    'net leggings',
    'net runner',
    'net tights short',
    'nets pants',
    'never give up leggings',
    'never hike alone',
    'never summer womens',
    'new balance leggings',
    'new balance leggings women',
    'new balance running capris women',
    'new clothing for women 2018',
    'new gear medical',
    'new order technique',
    'new year workout tank',
    'new years clothing for women',
    'night time muscle recovery',
    'night time running gear',
    'night time shorts for women',
    'nighttime running gear',
    'nighttime training pants',
    'nobo leggings',
    'noke leggings',
    'non chaffing stick',
    'norm core',
    'not your average runner',
    'noticing god',
    'now capri',
    'now i rise',
    'now leggings',
    'nylon calf compression',
    'nylon capri pants for women',
    'nylon leggings for women',
    'nylon spandex leggings',
    'nylon warmup pants',
    'nylon wedding rings for men',
    'oh la la exo',
    'oiselle pants',
    'oiselle running women',
    'oiselle tights',
    'old navy compression leggings women',
    'old rugby ball',
    'on joint support',
    'one ring of power',
    'opaque black tights for women',
    'opaque brown tights',
    'opaque support tights for women',
    'opaque thick tights',
    'opaque tight',
    'opaque tights black',
    'opaque tights brown',
    'opaque tights gray',
    'opaque tights grey',
    'opaque tights men',
    'opaque tights plus size',
    'opaque tights red',
    'opaque tights shiny',
    'opedix compression tights',
    'opedix tights',
    'opedix womens',
    'open knee leggings women',
    'outdoor capris for women',
    'outdoor master pants',
    'outdoor running tights',
    'outdoor survival clothing',
    'outdoor survival pants',
    'over garment for women',
    'over knee tights',
    'over the black rainbow',
    'over the foot yoga leggings',
    'over the foot yoga pants for women',
    'over the knee compression',
    'over the knee leggings',
    'overall capris',
    'overall capris for women',
    'overall for winter women',
    'overall leggings',
    'overalls for women black',
    'overalls for women insulated',
    'pace race',
    'pack of leggings women',
    'pack of tights for women',
    'pack tights women',
    'package leggings for women',
    'package tracking amazon',
    'pain figure',
    'pain for clothing',
    'pain in hip and thigh',
    'pain in hip area',
    'pain in right hip area',
    'pain in thigh and hip',
    'pain in thigh joint',
    'pain runner pro',
    'pain yoga',
    'panel leggings for women',
    'pant capri women',
    'pant for sport',
    'pant for woman',
    'pant hose',
    'pant leggings for women',
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
